$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force column D to Text format so numeric-looking strings
# (e.g. "233.73", "37.732.98") are stored as text, matching the source data.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "37.732.98"
$ws.Range("E2").Value = "  -0.22%  "

$ws.Range("D3").Value = "2.077.66"
$ws.Range("E3").Value = "  -1.26%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "233.73"
$ws.Range("E5").Value = "  -0.39%  "

$ws.Range("D6").Value = "0.623"
$ws.Range("E6").Value = "  -0.20%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").Value = "58.25"
$ws.Range("E8").Value = "  +0.40%  "

$ws.Range("E9").Value = "  +1.05%  "

$ws.Range("D10").Value = "0.0784"
$ws.Range("E10").Value = "  +0.74%  "

$ws.Range("E11").Value = "  +2.15%  "

$ws.Range("D12").Value = "2.384.77"
$ws.Range("E12").Value = "  -1.19%  "

$ws.Range("D13").Value = "14.77"
$ws.Range("E13").Value = "  +2.01%  "

$ws.Range("D14").Value = "20.84"
$ws.Range("E14").Value = "  -1.57%  "

$ws.Range("E15").Value = "  -1.36%  "

$ws.Range("E16").Value = "  +1.53%  "

$ws.Range("D17").Value = "2.068.83"
$ws.Range("E17").Value = "  -1.62%  "

$ws.Range("D18").Value = "37.700.90"
$ws.Range("E18").Value = "  -0.17%  "

$ws.Range("D19").Value = "6.19"
$ws.Range("E19").Value = "  -0.97%  "

$ws.Range("D20").Value = "71.09"
$ws.Range("E20").Value = "  +1.10%  "

$ws.Range("D21").Value = "0.0₃0833"
$ws.Range("E21").Value = "  +1.31%  "

$ws.Range("D22").Value = "228.51"
$ws.Range("E22").Value = "  +0.49%  "

$ws.Range("E23").Value = "  -0.03%  "

$ws.Range("D24").Value = "2.39"
$ws.Range("E24").Value = "  -0.94%  "

$ws.Range("E25").Value = "  +0.10%  "

$ws.Range("D26").Value = "170.63"
$ws.Range("E26").Value = "  +1.68%  "

$ws.Range("E27").Value = "  +5.36%  "

$ws.Range("D28").Value = "9.05"
$ws.Range("E28").Value = "  +1.12%  "

$ws.Range("D29").Value = "19.44"
$ws.Range("E29").Value = "  -0.34%  "

$ws.Range("E30").Value = "  -2.46%  "

$ws.Range("E31").Value = "  +2.47%  "

$ws.Range("E32").Value = "  +0.81%  "

$ws.Range("D33").Value = "0.0629"
$ws.Range("E33").Value = "  +1.15%  "

$ws.Range("E34").Value = "  +1.21%  "

$ws.Range("D35").Value = "2.49"
$ws.Range("E35").Value = "  -4.64%  "

$ws.Range("E36").Value = "  +2.20%  "

$ws.Range("E37").Value = "  -2.96%  "

$ws.Range("E38").Value = "  -0.24%  "

$ws.Range("D39").Value = "5.30"
$ws.Range("E39").Value = "  -2.76%  "

$ws.Range("D40").Value = "0.0967"
$ws.Range("E40").Value = "  -3.87%  "

$ws.Range("D41").Value = "98.13"
$ws.Range("E41").Value = "  +1.50%  "

$ws.Range("E42").Value = "  -2.31%  "

$ws.Range("D43").Value = "0.0214"
$ws.Range("E43").Value = "  +0.60%  "

$ws.Range("D44").Value = "1.449.46"
$ws.Range("E44").Value = "  -1.75%  "

$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").Value = "16.45"
$ws.Range("E45").Value = "  +6.84%  "

$ws.Range("E46").Value = "  -1.40%  "

$ws.Range("B47").Value = "FTXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D47").Value = "4.23"
$ws.Range("E47").Value = "  +2.41%  "

$ws.Range("E48").Value = "  +0.82%  "

$ws.Range("D49").Value = "7.42"
$ws.Range("E49").Value = "  +1.47%  "

$ws.Range("E50").Value = "  -0.72%  "

$ws.Range("D51").Value = "2.269.50"
$ws.Range("E51").Value = "  -1.34%  "

# Restore default (General) styling on column D so no stray number format lingers
$ws.Range("D2:D51").ClearFormats()
